# Remove the 'mode' column (column A) from every sheet that still has it.
# The TestSuite slot usage lets the content infer mode, so the explicit
# 'mode' enumerated field (Enumerated/Generated) is dropped; any other
# data validations shift left along with their columns automatically.

$wb = $excel.ActiveWorkbook

$sheetsWithMode = @(
    "TestCase",
    "TestSuite",
    "AcceptanceTestCase",
    "AcceptanceTestSuite",
    "StandardsComplianceTestSuite",
    "OneHopTestSuite"
)

foreach ($sheetName in $sheetsWithMode) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Columns("A").Delete()
}
